# Update cryptocurrency price/volume figures per refreshed source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $value) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "27.933.34"
Set-TextCell "E2" "  +0.44%  "
# Row 3
Set-TextCell "D3" "1.812.40"
Set-TextCell "E3" "  +1.70%  "
# Row 4
Set-TextCell "D4" "1.003"
Set-TextCell "E4" "  +0.01%  "
# Row 5
Set-TextCell "D5" "310.55"
Set-TextCell "E5" "  +0.05%  "
# Row 6
Set-TextCell "D6" "1.002"
Set-TextCell "E6" "  +0.02%  "
# Row 7
Set-TextCell "D7" "0.4958"
Set-TextCell "E7" "  -3.42%  "
# Row 8
Set-TextCell "D8" "0.3906"
Set-TextCell "E8" "  +2.98%  "
# Row 9
Set-TextCell "D9" "0.09701"
Set-TextCell "E9" "  +25.06%  "
# Row 10
Set-TextCell "E10" "  +1.46%  "
# Row 11
Set-TextCell "D11" "40.98"
Set-TextCell "E11" "  -0.39%  "
# Row 12
Set-TextCell "D12" "6.424"
Set-TextCell "E12" "  +3.64%  "
# Row 13
Set-TextCell "E13" "  +1.99%  "
# Row 14
Set-TextCell "D14" "1.002"
Set-TextCell "E14" "  +0.02%  "
# Row 15
Set-TextCell "D15" "1.812.77"
Set-TextCell "E15" "  +2.44%  "
# Row 16
Set-TextCell "D16" "7.280"
Set-TextCell "E16" "  +1.67%  "
# Row 17
Set-TextCell "E17" "  +5.52%  "
# Row 18
Set-TextCell "D18" "92.52"
Set-TextCell "E18" "  +1.21%  "
# Row 19
Set-TextCell "D19" "0.06675"
Set-TextCell "E19" "  +2.27%  "
# Row 20
Set-TextCell "E20" "  -0.01%  "
# Row 21
Set-TextCell "E21" "  +0.69%  "
# Row 22
Set-TextCell "D22" "5.912"
Set-TextCell "E22" "  +0.03%  "
# Row 23
Set-TextCell "D23" "27.985.37"
Set-TextCell "E23" "  +0.44%  "
# Row 24
Set-TextCell "D24" "11.16"
Set-TextCell "E24" "  +1.64%  "
# Row 25
Set-TextCell "E25" "  +0.27%  "
# Row 26
Set-TextCell "D26" "159.14"
Set-TextCell "E26" "  +0.30%  "
# Row 27
Set-TextCell "D27" "2.020.80"
Set-TextCell "E27" "  +2.02%  "
# Row 28
Set-TextCell "E28" "  +1.81%  "
# Row 29
Set-TextCell "D29" "2.394"
Set-TextCell "E29" "  +1.52%  "
# Row 30
Set-TextCell "D30" "128.04"
Set-TextCell "E30" "  +2.09%  "
# Row 31
Set-TextCell "D31" "0.1062"
Set-TextCell "E31" "  -1.30%  "
# Row 32
Set-TextCell "D32" "1.037"
Set-TextCell "E32" "  +1.04%  "
# Row 33
Set-TextCell "D33" "5.562"
Set-TextCell "E33" "  +1.53%  "
# Row 34
Set-TextCell "D34" "3.631"
Set-TextCell "E34" "  +0.66%  "
# Row 35
Set-TextCell "D35" "0.06697"
Set-TextCell "E35" "  -5.53%  "
# Row 36
Set-TextCell "D36" "8.946"
Set-TextCell "E36" "  +3.52%  "
# Row 38
Set-TextCell "D38" "0.2136"
Set-TextCell "E38" "  +0.74%  "
# Row 39
Set-TextCell "D39" "4.940"
# Row 40
Set-TextCell "D40" "11.25"
Set-TextCell "E40" "  -2.31%  "
# Row 41
Set-TextCell "D41" "0.6180"
Set-TextCell "E41" "  +1.58%  "
# Row 42
Set-TextCell "E42" "  +0.02%  "
# Row 43
Set-TextCell "D43" "1.146"
Set-TextCell "E43" "  -0.36%  "
# Row 44
Set-TextCell "D44" "13.10"
Set-TextCell "E44" "  +0.27%  "
# Row 45
Set-TextCell "D45" "0.5874"
Set-TextCell "E45" "  -1.44%  "
# Row 46
Set-TextCell "D46" "3.691"
Set-TextCell "E46" "  -0.48%  "
# Row 47
Set-TextCell "E47" "  -2.87%  "
# Row 48
Set-TextCell "D48" "122.84"
Set-TextCell "E48" "  -3.15%  "
# Row 49
Set-TextCell "D49" "1.936"
Set-TextCell "E49" "  +2.13%  "
# Row 50
Set-TextCell "D50" "1.178"
Set-TextCell "E50" "  -2.97%  "
# Row 51
Set-TextCell "D51" "0.06793"
Set-TextCell "E51" "  +1.29%  "
